# Adds a 'hole_id' index column (A1:A32) to the "train" worksheet so that
# cross validation can line samples up by hole id.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$holeIds = @(
    "BRG_13_01",
    "BRG_05_13",
    "BRG_05_10",
    "BRG_05_14",
    "ECO_09_04",
    "BRG_05_05",
    "BRG_05_02",
    "BRG_16_08",
    "BRG_05_15",
    "BRG_16_07",
    "BRG_16_09",
    "BRG_05_12",
    "ECO_09_02",
    "BRG_16_04B",
    "BRG_05_09",
    "BRG_01_03",
    "BRG_16_01",
    "BRG_16_03",
    "BRG_01_07",
    "ECO_09_01",
    "BRG_16_05",
    "BRG_08_01",
    "ECO_09_03",
    "BRG_01_01",
    "BRG_01_04",
    "BRG_16_06",
    "BRG_05_01",
    "BRG_01_09",
    "BRG_13_02",
    "BRG_01_02",
    "BRG_01_06"
)

# Header cell for the new index column.
$ws.Range("A1").Value = "hole_id"

# Match the bold/centered/bordered header style already used by B1:M1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# Replace the old plain numeric row index (0..30) with the hole id text.
for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
